# Update cryptocurrency price/volume data to the latest scrape.
# Every written cell is forced to Text (General-looking numbers like
# "1.005" would otherwise auto-coerce to a Number on assignment), then
# ClearFormats() drops the temporary "@" number-format style so the cell
# keeps using the sheet default style (matches the original inlineStr cells).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = '@'
    $c.Value = $text
    $c.ClearFormats()
}

Set-TextValue 'D2' '28.197.85'
Set-TextValue 'E2' '  +0.86%  '
Set-TextValue 'D3' '1.882.62'
Set-TextValue 'E3' '  +1.25%  '
Set-TextValue 'D4' '1.005'
Set-TextValue 'E4' '  +0.22%  '
Set-TextValue 'D5' '313.82'
Set-TextValue 'E5' '  +0.74%  '
Set-TextValue 'D6' '1.005'
Set-TextValue 'E6' '  +0.27%  '
Set-TextValue 'D7' '0.5141'
Set-TextValue 'E7' '  +0.09%  '
Set-TextValue 'D8' '0.3919'
Set-TextValue 'E8' '  +3.00%  '
Set-TextValue 'D9' '0.08362'
Set-TextValue 'E9' '  +1.22%  '
Set-TextValue 'D10' '1.123'
Set-TextValue 'E10' '  +1.42%  '
Set-TextValue 'D11' '41.67'
Set-TextValue 'E11' '  +0.15%  '
Set-TextValue 'D12' '6.226'
Set-TextValue 'E12' '  +0.73%  '
Set-TextValue 'B13' 'WrappedEther'
Set-TextValue 'C13' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 'D13' '1.887.06'
Set-TextValue 'E13' '  +1.53%  '
Set-TextValue 'B14' 'Solana'
Set-TextValue 'C14' 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextValue 'D14' '20.70'
Set-TextValue 'E14' '  +1.34%  '
Set-TextValue 'D15' '7.275'
Set-TextValue 'E15' '  +1.12%  '
Set-TextValue 'D16' '1.006'
Set-TextValue 'E16' '  +0.30%  '
Set-TextValue 'D17' '0.00001105'
Set-TextValue 'E17' '  +1.20%  '
Set-TextValue 'D18' '91.28'
Set-TextValue 'E18' '  +1.09%  '
Set-TextValue 'D19' '0.06663'
Set-TextValue 'E19' '  +0.92%  '
Set-TextValue 'D20' '17.79'
Set-TextValue 'E20' '  +0.38%  '
Set-TextValue 'D21' '1.005'
Set-TextValue 'E21' '  +0.35%  '
Set-TextValue 'D22' '6.052'
Set-TextValue 'E22' '  +0.96%  '
Set-TextValue 'D23' '28.254.36'
Set-TextValue 'E23' '  +0.94%  '
Set-TextValue 'D24' '11.17'
Set-TextValue 'E24' '  +1.55%  '
Set-TextValue 'E25' '  +2.43%  '
Set-TextValue 'D26' '2.091.71'
Set-TextValue 'E26' '  +0.83%  '
Set-TextValue 'E27' '  -2.77%  '
Set-TextValue 'D28' '159.69'
Set-TextValue 'E29' '  +1.51%  '
Set-TextValue 'D30' '125.32'
Set-TextValue 'E30' '  +0.78%  '
Set-TextValue 'D31' '0.1064'
Set-TextValue 'E31' '  +0.27%  '
Set-TextValue 'D32' '1.042'
Set-TextValue 'E32' '  +0.79%  '
Set-TextValue 'D33' '5.856'
Set-TextValue 'E33' '  +4.85%  '
Set-TextValue 'E34' '  +0.05%  '
Set-TextValue 'D35' '9.689'
Set-TextValue 'E35' '  +1.55%  '
Set-TextValue 'D36' '0.02465'
Set-TextValue 'E36' '  +1.74%  '
Set-TextValue 'D37' '0.06568'
Set-TextValue 'E37' '  +0.65%  '
Set-TextValue 'D38' '0.2189'
Set-TextValue 'E38' '  +0.43%  '
Set-TextValue 'D39' '1.207'
Set-TextValue 'E39' '  +0.17%  '
Set-TextValue 'D40' '0.6519'
Set-TextValue 'E40' '  +1.69%  '
Set-TextValue 'D41' '1.238'
Set-TextValue 'E41' '  +0.39%  '
Set-TextValue 'D42' '5.003'
Set-TextValue 'E42' '  +2.77%  '
Set-TextValue 'E43' '  +0.75%  '
Set-TextValue 'D44' '0.6149'
Set-TextValue 'E44' '  +0.86%  '
Set-TextValue 'D45' '13.17'
Set-TextValue 'E45' '  +0.99%  '
Set-TextValue 'D46' '1.289'
Set-TextValue 'E46' '  +0.60%  '
Set-TextValue 'E47' '  +0.97%  '
Set-TextValue 'D48' '2.022'
Set-TextValue 'E48' '  +2.61%  '
Set-TextValue 'D49' '1.236'
Set-TextValue 'E49' '  +2.70%  '
Set-TextValue 'D50' '121.04'
Set-TextValue 'E50' '  +0.42%  '
Set-TextValue 'D51' '78.73'
Set-TextValue 'E51' '  -0.91%  '
